# Updates the cryptocurrency price/volume table:
#  - Refreshes Price (D) and Volume(1h) (E) columns with latest scraped values
#  - Swaps the OKB / Stacks rows (38-39) including Coin name, Link, Price, Volume

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value renders as a plain number (e.g. '512.36') need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values instead of keeping the original text representation used in the sheet.
$textCells = @("D5", "D6", "D10", "D19", "D21", "D22", "D23", "D24", "D28", "D34", "D37", "D38", "D39", "D40", "D47", "D48", "D49")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '56.560.90'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '2.323.55'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '512.36'
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("D6").Value = '131.67'
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("E9").Value = '  -3.92%  '
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("E12").Value = '  -2.02%  '
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").Value = '2.737.62'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '56.524.83'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '2.326.95'
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").Value = '327.92'
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '61.42'
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("D24").Value = '8.60'
$ws.Range("E24").Value = '  +8.58%  '
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = '167.89'
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("D30").Value = '0.0₃0719'
$ws.Range("E30").Value = '  -4.05%  '
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("E36").Value = '  -2.87%  '
$ws.Range("D37").Value = '0.886'
$ws.Range("E37").Value = '  -4.40%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '1.56'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '38.56'
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("D40").Value = '148.66'
$ws.Range("E40").Value = '  +7.67%  '
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("E44").Value = '  -4.30%  '
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("E46").Value = '  -2.20%  '
$ws.Range("D47").Value = '0.556'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '18.19'
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").Value = '0.380'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  -1.58%  '
$ws.Range("E51").Value = '  +0.91%  '
